# Fill in the "checkboard" keyboard grid, extending it from 5 rows to 9 rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Formatting first: extend the existing "Dane wejściowe" / "Uwaga"
# styles down into the new rows (copy formats only, so we don't disturb
# values and don't create duplicate style records). ---
$ws.Range("B5:K5").Copy() | Out-Null
$ws.Range("B6:K9").PasteSpecial(-4122) | Out-Null

$ws.Range("A5").Copy() | Out-Null
$ws.Range("A6:A9").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Row 2: the first keyboard row now only keeps Q/W/E (rest cleared) ---
$ws.Range("E2").Value = $null
$ws.Range("F2").Value = $null
$ws.Range("G2").Value = $null
$ws.Range("H2").Value = $null
$ws.Range("I2").Value = $null
$ws.Range("J2").Value = $null
$ws.Range("K2").Value = $null

# --- Row 3 ---
$ws.Range("B3").Value = "R"
$ws.Range("C3").Value = "T"
$ws.Range("D3").Value = "Y"
$ws.Range("E3").Value = "U"
$ws.Range("F3").Value = "I"
$ws.Range("G3").Value = "O"
$ws.Range("H3").Value = "P"
$ws.Range("I3").Value = "A"
$ws.Range("J3").Value = "S"
$ws.Range("K3").Value = "D"

# --- Row 4 ---
$ws.Range("B4").Value = "F"
$ws.Range("C4").Value = "G"
$ws.Range("D4").Value = "H"
$ws.Range("E4").Value = "J"
$ws.Range("F4").Value = "K"
$ws.Range("G4").Value = "L"
$ws.Range("H4").Value = "Z"
$ws.Range("I4").Value = "X"
$ws.Range("J4").Value = "C"
$ws.Range("K4").Value = "V"

# --- Row 5 ---
$ws.Range("B5").Value = "B"
$ws.Range("C5").Value = "N"
$ws.Range("D5").Value = "M"
$ws.Range("E5").Value = "Ą"
$ws.Range("F5").Value = "Ć"
$ws.Range("G5").Value = "Ł"
$ws.Range("H5").Value = "Ń"
$ws.Range("I5").Value = "Ó"
$ws.Range("J5").Value = "Ż"
$ws.Range("K5").Value = "Ź"

# --- Row 6 (new) ---
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = ","
$ws.Range("C6").Value = "."
$ws.Range("D6").Value = "q"
$ws.Range("E6").Value = "w"
$ws.Range("F6").Value = "e"
$ws.Range("G6").Value = "r"
$ws.Range("H6").Value = "t"
$ws.Range("I6").Value = "y"
$ws.Range("J6").Value = "u"
$ws.Range("K6").Value = "i"

# --- Row 7 (new) ---
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "o"
$ws.Range("C7").Value = "p"
$ws.Range("D7").Value = "a"
$ws.Range("E7").Value = "s"
$ws.Range("F7").Value = "d"
$ws.Range("G7").Value = "f"
$ws.Range("H7").Value = "g"
$ws.Range("I7").Value = "h"
$ws.Range("J7").Value = "j"
$ws.Range("K7").Value = "k"

# --- Row 8 (new) ---
$ws.Range("A8").Value = 9
$ws.Range("B8").Value = "l"
$ws.Range("C8").Value = "z"
$ws.Range("D8").Value = "x"
$ws.Range("E8").Value = "c"
$ws.Range("F8").Value = "v"
$ws.Range("G8").Value = "b"
$ws.Range("H8").Value = "n"
$ws.Range("I8").Value = "m"
$ws.Range("J8").Value = "ą"
$ws.Range("K8").Value = "ć"

# --- Row 9 (new) ---
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "ł"
$ws.Range("C9").Value = "ń"
$ws.Range("D9").Value = "ó"
$ws.Range("E9").Value = "ż"
$ws.Range("F9").Value = "ź"
$ws.Range("G9").Value = "ś"
$ws.Range("H9").Value = "Ś"

# --- Column B gets an explicit width, matching the rest of the grid columns ---
$ws.Columns("B").ColumnWidth = 2.8

# --- Selection now sits on the last-filled cell of the grid ---
$ws.Range("H9").Select()

$wb.Save()
